$d = $word.ActiveDocument
$replacements = @(
    ,@("9-0=9", "30+28=58")
    ,@("91+7=98", "80-29=51")
    ,@("55+16=71", "55-33=22")
    ,@("68+16=84", "14+73=87")
    ,@("65+17=82", "29-13=16")
    ,@("72-58=14", "30+54=84")
    ,@("52+35=87", "91-43=48")
    ,@("74-23=51", "45-3=42")
    ,@("55+6=61", "76-24=52")
    ,@("44+32=76", "98-39=59")
    ,@("95-37=58", "46+27=73")
    ,@("4+26=30", "59+36=95")
    ,@("45-17=28", "67-53=14")
    ,@("87-81=6", "51-27=24")
    ,@("49-37=12", "5+53=58")
    ,@("25+5=30", "38+46=84")
    ,@("60+25=85", "56-15=41")
    ,@("22+8=30", "24+62=86")
    ,@("72+12=84", "77-32=45")
    ,@("28+40=68", "93-14=79")
    ,@("32-1=31", "50-25=25")
    ,@("41-4=37", "96-57=39")
    ,@("12+87=99", "82-61=21")
    ,@("53-12=41", "84+1=85")
    ,@("65-59=6", "89-8=81")
    ,@("84-13=71", "64-5=59")
    ,@("2+72=74", "98-26=72")
    ,@("97-57=40", "44+41=85")
    ,@("72+4=76", "89-77=12")
    ,@("80-24=56", "53+36=89")
    ,@("19+25=44", "73-26=47")
    ,@("0+49=49", "20+70=90")
    ,@("2+11=13", "45+29=74")
    ,@("44+1=45", "68-38=30")
    ,@("73+15=88", "77-19=58")
    ,@("85-69=16", "35+49=84")
    ,@("45+45=90", "11+33=44")
    ,@("23+41=64", "2+9=11")
    ,@("15+24=39", "48-47=1")
    ,@("3+22=25", "12+20=32")
    ,@("59-8=51", "90-16=74")
    ,@("49+40=89", "61-34=27")
    ,@("70+20=90", "84+6=90")
    ,@("82-41=41", "23-10=13")
    ,@("75+8=83", "11+78=89")
    ,@("15+13=28", "49-4=45")
    ,@("84+3=87", "42+40=82")
    ,@("37-32=5", "84-30=54")
    ,@("83-3=80", "96-9=87")
    ,@("72-65=7", "29+31=60")
    ,@("77-52=25", "49-13=36")
    ,@("0+25=25", "52+15=67")
    ,@("8+61=69", "67+21=88")
    ,@("84+0=84", "25+70=95")
    ,@("58-53=5", "43-13=30")
    ,@("23+61=84", "25+62=87")
    ,@("4-4=0", "86-13=73")
    ,@("46-0=46", "52+28=80")
    ,@("76-73=3", "77-26=51")
    ,@("31+40=71", "77-32=45")
    ,@("39+9=48", "14+48=62")
    ,@("66-3=63", "21-15=6")
    ,@("67-26=41", "6+1=7")
    ,@("49+4=53", "10+26=36")
    ,@("47-46=1", "45-14=31")
    ,@("89-4=85", "70-25=45")
    ,@("37+20=57", "2+8=10")
    ,@("11+1=12", "44-42=2")
    ,@("77+5=82", "10+12=22")
    ,@("22+62=84", "91-16=75")
    ,@("44+44=88", "48-9=39")
    ,@("30+41=71", "14+1=15")
    ,@("8+59=67", "87-36=51")
    ,@("3+13=16", "85-79=6")
    ,@("47+14=61", "34-25=9")
    ,@("45-5=40", "99-5=94")
    ,@("26+16=42", "19+14=33")
    ,@("7+2=9", "47-42=5")
    ,@("68-52=16", "23+15=38")
    ,@("81-32=49", "13+38=51")
    ,@("84-47=37", "66-2=64")
    ,@("21+39=60", "12-4=8")
    ,@("62+0=62", "91-40=51")
    ,@("29+67=96", "17+75=92")
    ,@("54-2=52", "87-87=0")
    ,@("30+13=43", "10+45=55")
    ,@("67+4=71", "90-60=30")
    ,@("42+10=52", "1+68=69")
    ,@("39+23=62", "58-31=27")
    ,@("31+5=36", "99-22=77")
    ,@("78-51=27", "49-12=37")
    ,@("43-10=33", "64-61=3")
    ,@("14+84=98", "61-50=11")
    ,@("38+49=87", "44-35=9")
    ,@("26+46=72", "35-13=22")
    ,@("83-61=22", "82+15=97")
    ,@("68-21=47", "79-1=78")
    ,@("60+30=90", "24-6=18")
    ,@("84+5=89", "23+57=80")
    ,@("57-7=50", "82-50=32")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}
Write-Host "Done applying replacements"
